$d = $word.ActiveDocument

# The CV header currently renders as "Ben Barrrrr": the last name "Barr" is
# followed by a stray extra run of "rrrr" that shouldn't be there. Remove
# just that trailing stray run so the header reads "Ben Barr" again, while
# leaving every other run (including the "r" run right before it) intact.

$strayText = "rrrr"

$anchor = $d.Content
$found = $anchor.Find.Execute("Ben Bar", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $target = $d.Range($anchor.End, $anchor.End + $strayText.Length)
    if ($target.Text -eq $strayText) {
        $target.Delete()
    }
}
